# Ajuste esqueleto CS_06_02_CO y otros...
# Correccion de mayusculas en titulos de recursos aprovechados
#
# The "CUADERNO DE ESTUDIO" sheet lists, in column H, the titles of
# resources that were "aprovechados" (leveraged/reused). Several of those
# titles were typed with a leading capital letter that didn't match the
# canonical (lower-case, sentence-style) title used elsewhere in the
# workbook. This fixes the casing (and one missing accent) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUADERNO DE ESTUDIO")

$ws.Range("H15").Value = "El paleolítico"
$ws.Range("H17").Value = "Identifica las características del paleolítico"
$ws.Range("H26").Value = "El neolítico"
$ws.Range("H31").Value = "Refuerza el aprendizaje: El neolítico"
$ws.Range("H35").Value = "El neolítico y la edad de los metales"
$ws.Range("H37").Value = "Relaciona conceptos con su período histórico"
$ws.Range("H43").Value = "Arte y religión en la edad de piedra"
$ws.Range("H52").Value = "Refuerza tu aprendizaje: El arte en la prehistoria"
$ws.Range("H54").Value = "Competencias: comentario de una obra de arte de la prehistoria"

# Restore the selection/scroll position that Excel recorded when the file
# was last saved.
$ws.Activate() | Out-Null
$ws.Range("H25").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 7
